$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# A new log entry was captured; insert a fresh row above row 29 which
# pushes every existing row (29-50) down by one (A50 "Broadband" -> A51,
# the September log rows R/S -> next row, the August P/Q rows -> next row,
# etc.), then populate the new row 29 with the new entry.
$ws.Rows.Item(29).Insert()

$ws.Cells.Item(29, 18).Value = "logging iob internet"
$ws.Cells.Item(29, 19).Value = "2024-09-03 20:09:12"
